# Fill in the four newly-solved Linked List problems (rows 87-90) on Sheet1.
# Columns: A=ID, B=Topic, C=Problem Name, D=Difficulty, E=Status,
#          F=Date Solved, G=Time Complexity, H=Space Complexity,
#          I=Approach Summary

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use an already-formatted "Date Solved" cell as the format donor so the
# new date cells reuse the existing date style instead of minting a new one.
$dateDonor = $ws.Range("F86")

$rows = @(
    @{ Row = 87; Topic = "Linked List"; Problem = "Linked List Cycle II";          Difficulty = "Medium"; Status = "Done"; Date = 45919; Time = "O(n)";     Space = "O(1)"; Approach = "Slow & Fast Pointers" },
    @{ Row = 88; Topic = "Linked List"; Problem = "Intersection of 2 Linked List"; Difficulty = "Medium"; Status = "Done"; Date = 45919; Time = "O(m + n)"; Space = "O(1)"; Approach = $null },
    @{ Row = 89; Topic = "Linked List"; Problem = "Remove Linked List Elements";   Difficulty = "Medium"; Status = "Done"; Date = 45919; Time = "O(n)";     Space = "O(1)"; Approach = $null },
    @{ Row = 90; Topic = "Linked List"; Problem = "Delete Node in a Linked List";  Difficulty = "Medium"; Status = "Done"; Date = 45919; Time = "O(1)";     Space = "O(1)"; Approach = $null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Range("B$rowNum").Value = $r.Topic
    $ws.Range("C$rowNum").Value = $r.Problem
    $ws.Range("D$rowNum").Value = $r.Difficulty
    $ws.Range("E$rowNum").Value = $r.Status

    $dateDonor.Copy()
    $ws.Range("F$rowNum").PasteSpecial(-4122)
    $ws.Range("F$rowNum").Value = $r.Date

    $ws.Range("G$rowNum").Value = $r.Time
    $ws.Range("H$rowNum").Value = $r.Space

    if ($r.Approach) {
        $ws.Range("I$rowNum").Value = $r.Approach
    }
}

$excel.CutCopyMode = 0
$ws.Range("B91").Select() | Out-Null
